# Fruta / hortaliza, semanal
# Re-order the daily price rows (2-5) into their weekly-sorted positions.
# Row 2 <-> Row 5 swap fully (date, volume, price range, unit, origin, price/kg, units).
# Row 3 <-> Row 4 swap their date and volume values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 becomes what was Row 5 ---
$ws.Range("D2").Value = 44687
$ws.Range("J2").Value = 160
$ws.Range("K2").Value = 3000
$ws.Range("L2").Value = 3500
$ws.Range("M2").Value = 3250
$ws.Range("N2").Value = "$/docena de matas"
$ws.Range("O2").Value = "Región Metropolitana"
$ws.Range("P2").Value = 542
$ws.Range("Q2").Value = 6

# --- Row 3 keeps same K/L/M/N/O/P/Q, only date/volume change ---
$ws.Range("D3").Value = 44691
$ws.Range("J3").Value = 100

# --- Row 4 keeps same K/L/M/N/O/P/Q, only date/volume change ---
$ws.Range("D4").Value = 44692
$ws.Range("J4").Value = 120

# --- Row 5 becomes what was Row 2 ---
$ws.Range("D5").Value = 44221
$ws.Range("J5").Value = 250
$ws.Range("K5").Value = 1300
$ws.Range("L5").Value = 1500
$ws.Range("M5").Value = 1420
$ws.Range("N5").Value = "$/atado"
$ws.Range("O5").Value = "Provincia de Diguillín"
$ws.Range("P5").Value = 1420
$ws.Range("Q5").Value = 1
